# Apply the "Add files via upload" edit:
#   - View_Print!A1:O1 loses its static shared-string header row; A1 becomes
#     a formula mirroring A2 (=A2), and B1:O1 are cleared (styles kept).
#   - View_Print selection moves from A1:O1 (active O1) to B2:P2 (active B2).
#   - Input!A2 gets the value 1.
#   - Input selection moves from A2 to A3.
# The Input sheet is touched last so it ends up the active/selected tab,
# matching the saved workbook (tabSelected stays on Input, not View_Print).

$wb = $excel.ActiveWorkbook

$wsPrint = $wb.Worksheets.Item("View_Print")
$wsPrint.Range("A1").Formula = "=A2"
$wsPrint.Range("B1:O1").ClearContents()
$wsPrint.Range("B2:P2").Select() | Out-Null

$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("A2").Value = 1
$wsInput.Range("A3").Select() | Out-Null
